# FormulaTemplate.xlsx update:
#  - Add three new worksheets: "Copy Right", "ReplaceTest", "Outside Reference"
#  - Populate them with JETT formula-test content
#  - Add a new row (row 6) to the "Formula Test" sheet

$wb = $excel.ActiveWorkbook

# --- Remember / anchor on the original active sheet ---
$formulaTest = $wb.Worksheets.Item("Formula Test")

# --- Add "Copy Right" sheet after the last existing sheet (MultiLevel) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsCopyRight = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$wsCopyRight.Name = "Copy Right"

$wsCopyRight.Range("A1").Value = '<jt:for start="1" end="10" var="n" copyRight="true">${n}'
$wsCopyRight.Range("A2").Value = '${2*n}'
$wsCopyRight.Range("A3").Value = '$[SUM(A1+A2)]'
$wsCopyRight.Range("A4").Value = '</jt:for>'
$wsCopyRight.Range("A1:A4").Select() | Out-Null

# --- Add "ReplaceTest" sheet after "Copy Right" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsReplaceTest = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$wsReplaceTest.Name = "ReplaceTest"

$wsReplaceTest.Range("A1").Value = '<jt:for start="1" end="10" var="n">${n}'
$wsReplaceTest.Range("A2").Value = '${2*n}'
$wsReplaceTest.Range("A3").Value = '$[SUM(A1+A2)]'
$wsReplaceTest.Range("A4").Value = '</jt:for>'

# --- Add "Outside Reference" sheet after "ReplaceTest" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsOutsideRef = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$wsOutsideRef.Name = "Outside Reference"

$wsOutsideRef.Range("A1").Value = '${two}'
$wsOutsideRef.Range("B1").Value = '<jt:forEach items="${primes}" var="x">${x}'
$wsOutsideRef.Range("C1").Value = '$[A1 * B1]'
$wsOutsideRef.Range("D1").Value = '<jt:forEach items="${morePrimes}" var="y">${y}'
$wsOutsideRef.Range("E1").Value = '$[A1 * B1 * D1]'
$wsOutsideRef.Range("F1").Value = '</jt:forEach></jt:forEach>'

# --- Add new row 6 to the "Formula Test" sheet ---
$formulaTest.Range("A6").Value = "Population Different?"
$formulaTest.Range("C6").Value = '$[B4 <> H4]'

# --- Restore the originally active sheet & selection ---
$formulaTest.Activate()
$formulaTest.Range("A1:F1").Select() | Out-Null
